$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.786.61"
$ws.Range("E2").Value = "  +2.68%  "

$ws.Range("D3").Value = "2.091.38"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  -0.04%  "

$cell = $ws.Range("D5")
$cell.Value = "'228.92"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("E6").Value = "  +0.96%  "

$cell = $ws.Range("D7")
$cell.Value = "'60.65"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "

$ws.Range("E8").Value = "  -0.07%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.385"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.60%  "

$cell = $ws.Range("D10")
$cell.Value = "'0.0837"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "2.401.44"
$ws.Range("E12").Value = "  +2.35%  "

$cell = $ws.Range("D13")
$cell.Value = "'15.00"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.61%  "

$cell = $ws.Range("D14")
$cell.Value = "'21.84"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.01%  "

$cell = $ws.Range("D15")
$cell.Value = "'0.798"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +4.57%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").Value = "2.089.18"
$ws.Range("E17").Value = "  +3.04%  "

$ws.Range("D18").Value = "38.675.92"
$ws.Range("E18").Value = "  +2.52%  "

$cell = $ws.Range("D19")
$cell.Value = "'71.81"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "

$cell = $ws.Range("D20")
$cell.Value = "'6.05"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +1.49%  "

$cell = $ws.Range("D22")
$cell.Value = "'227.02"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("E24").Value = "  -0.19%  "

$cell = $ws.Range("D25")
$cell.Value = "'2.35"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.37%  "

$cell = $ws.Range("D26")
$cell.Value = "'170.93"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "

$cell = $ws.Range("D27")
$cell.Value = "'9.47"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("E28").Value = "  +8.33%  "

$ws.Range("E29").Value = "  +12.95%  "

$ws.Range("E30").Value = "  +2.24%  "

$ws.Range("E31").Value = "  +1.07%  "

$cell = $ws.Range("D32")
$cell.Value = "'2.39"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +6.11%  "

$ws.Range("E33").Value = "  +3.08%  "

$ws.Range("E34").Value = "  +4.81%  "

$cell = $ws.Range("D35")
$cell.Value = "'0.0613"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "

$ws.Range("E36").Value = "  +1.99%  "

$cell = $ws.Range("D37")
$cell.Value = "'6.42"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("E38").Value = "  +3.33%  "

$ws.Range("E39").Value = "  +0.01%  "

$cell = $ws.Range("D40")
$cell.Value = "'18.22"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "1.541.96"
$ws.Range("E41").Value = "  +0.45%  "

$cell = $ws.Range("D42")
$cell.Value = "'100.79"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("E43").Value = "  +4.36%  "

$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("E46").Value = "  +9.06%  "

$ws.Range("E47").Value = "  +1.34%  "

$cell = $ws.Range("D48")
$cell.Value = "'4.10"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("E49").Value = "  +2.76%  "

$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").Value = "2.287.51"
$ws.Range("E51").Value = "  +2.39%  "
